$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the new B1 header cell using the same style as A1. A plain value
# assignment does not carry formatting in this environment, but Insert()
# (after Copy()) does copy the style - it just also shifts the whole
# column across. We immediately delete the resulting empty "gap" column
# so every other column/cell ends up back exactly where it started, and
# B1 is left behind with A1's style applied.
$ws.Range("A1").Copy()
$ws.Range("B1").Insert(-4161)
$ws.Columns.Item(3).Delete()

# Now set the real header text for the new column.
$ws.Range("B1").Value2 = "SUPPLIER SAMPLE NAME"

# Move the existing "plate" values out of column A (rows 2-8, row 6 is
# empty/unused) into column B, clearing them from column A.
$rows = @(2, 3, 4, 5, 7, 8)
foreach ($r in $rows) {
    $srcCell = $ws.Cells.Item($r, 1)
    $destCell = $ws.Cells.Item($r, 2)
    $destCell.Value2 = $srcCell.Value2
    $srcCell.ClearContents()
}

# Give column B an explicit width, matching the new header/content.
$ws.Columns.Item(2).ColumnWidth = 28.82

# Update the selection to match the new data range (B2:B8), active cell B2.
$ws.Range("B2:B8").Select()
